$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Insert a brand-new row (id=1) right above the row that currently
#     holds id=2 (the first data row after the header row). ---
$targetRow = $t.Rows.Item(2)
$newRow = $t.Rows.Add($targetRow)

# Rows.Add() only materializes a single w:tc in the new row in this
# runtime, so grow it back out to the table's real column count by
# repeatedly splitting the first cell in two.
while ($t.Rows.Item(2).Cells.Count -lt $t.Columns.Count) {
    $t.Cell(2, 1).Split(1, 2)
}

$t.Cell(2, 1).Range.Text = "1"
$t.Cell(2, 2).Range.Text = "2023-04-12"
$t.Cell(2, 3).Range.Text = "собирается"
$t.Cell(2, 4).Range.Text = "25.22"

# --- Update the weight values for the id=2 and id=5 rows. ---
$d.Content.Find.Execute("20.0", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "21.04", 2)
$d.Content.Find.Execute("24.0", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "25.36", 2)
